# Auto-generated edit script: updates market-price derived columns (H-N)
# on sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR to reflect refreshed
# Universalis market data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 3307.0527
$ws.Range("I28").Value = 1206.6
$ws.Range("J28").Value = 5640.8887
$ws.Range("K28").Value = 1206.6
$ws.Range("L28").Value = 5640.8887
$ws.Range("M28").Value = -721.5999999999999
$ws.Range("N28").Value = -6610.8887
$ws.Range("H33").Value = 74.375
$ws.Range("I33").Value = 24.166666
$ws.Range("K33").Value = 24.166666
$ws.Range("M33").Value = 204.833334
$ws.Range("H53").Value = 458.1
$ws.Range("I53").Value = 550.1429000000001
$ws.Range("K53").Value = 550.1429000000001
$ws.Range("M53").Value = 86.85709999999995
$ws.Range("H62").Value = 5064.143
$ws.Range("I62").Value = 1877
$ws.Range("K62").Value = 1877
$ws.Range("M62").Value = -1253
$ws.Range("H65").Value = 5064.143
$ws.Range("I65").Value = 1877
$ws.Range("K65").Value = 9385
$ws.Range("M65").Value = -6265
$ws.Range("H86").Value = 10000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 10000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 10000
$ws.Range("N86").Value = -12246
$ws.Range("H89").Value = 10000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 10000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 50000
$ws.Range("N89").Value = -61232
$ws.Range("H96").Value = 4598.5
$ws.Range("I96").Value = 4598
$ws.Range("J96").Value = 4599
$ws.Range("K96").Value = 13794
$ws.Range("L96").Value = 13797
$ws.Range("M96").Value = -12421
$ws.Range("N96").Value = -16543
$ws.Range("H98").Value = 333.33334
$ws.Range("I98").Value = 333.33334
$ws.Range("K98").Value = 333.33334
$ws.Range("M98").Value = 1164.66666
$ws.Range("H100").Value = 306.5
$ws.Range("I100").Value = 306.5
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 306.5
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = 234.5
$ws.Range("H106").Value = 6118.3
$ws.Range("I106").Value = 6118.3
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 6118.3
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -5487.3
$ws.Range("H113").Value = 5771.6665
$ws.Range("I113").Value = 5771.6665
$ws.Range("K113").Value = 5771.6665
$ws.Range("M113").Value = -2517.6665
$ws.Range("H122").Value = 333.33334
$ws.Range("I122").Value = 333.33334
$ws.Range("K122").Value = 1000.00002
$ws.Range("M122").Value = 1449.99998
$ws.Range("H137").Value = 2586.3125
$ws.Range("I137").Value = 1896.4
$ws.Range("J137").Value = 2899.9092
$ws.Range("K137").Value = 5689.200000000001
$ws.Range("L137").Value = 8699.7276
$ws.Range("M137").Value = -3139.200000000001
$ws.Range("N137").Value = -13799.7276
$ws.Range("M86").ClearContents()
$ws.Range("M89").ClearContents()
$ws.Range("N100").ClearContents()
$ws.Range("N106").ClearContents()

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 934.75
$ws.Range("I2").Value = 934.75
$ws.Range("K2").Value = 934.75
$ws.Range("M2").Value = -821.75
$ws.Range("H61").Value = 7192.9
$ws.Range("I61").Value = 4732.25
$ws.Range("K61").Value = 4732.25
$ws.Range("M61").Value = -4520.25
$ws.Range("H102").Value = 3013.2104
$ws.Range("I102").Value = 1432.2858
$ws.Range("K102").Value = 1432.2858
$ws.Range("M102").Value = 189.7141999999999
$ws.Range("H110").Value = 499.6
$ws.Range("J110").Value = 174
$ws.Range("L110").Value = 174
$ws.Range("N110").Value = -4264
$ws.Range("H116").Value = 934.75
$ws.Range("I116").Value = 934.75
$ws.Range("K116").Value = 934.75
$ws.Range("M116").Value = 1359.25
$ws.Range("H132").Value = 2314.3
$ws.Range("I132").Value = 2314.3
$ws.Range("K132").Value = 6942.900000000001
$ws.Range("M132").Value = -4412.900000000001
$ws.Range("H136").Value = 7192.9
$ws.Range("I136").Value = 4732.25
$ws.Range("K136").Value = 14196.75
$ws.Range("M136").Value = -11646.75

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 934.75
$ws.Range("I3").Value = 934.75
$ws.Range("K3").Value = 934.75
$ws.Range("M3").Value = -820.75
$ws.Range("H20").Value = 3000
$ws.Range("I20").Value = 2500
$ws.Range("J20").Value = 3500
$ws.Range("K20").Value = 2500
$ws.Range("L20").Value = 3500
$ws.Range("M20").Value = -2253
$ws.Range("N20").Value = -3994
$ws.Range("H86").Value = 4158.3125
$ws.Range("J86").Value = 6535.5713
$ws.Range("L86").Value = 6535.5713
$ws.Range("N86").Value = -8781.5713
$ws.Range("H89").Value = 4158.3125
$ws.Range("J89").Value = 6535.5713
$ws.Range("L89").Value = 32677.8565
$ws.Range("N89").Value = -43909.85649999999
$ws.Range("H94").Value = 535.3077
$ws.Range("I94").Value = 450.9
$ws.Range("K94").Value = 450.9
$ws.Range("M94").Value = 0.1000000000000227
$ws.Range("H99").Value = 2453.6
$ws.Range("I99").Value = 2067.25
$ws.Range("K99").Value = 2067.25
$ws.Range("M99").Value = -569.25
$ws.Range("H134").Value = 1562.3
$ws.Range("I134").Value = 1732.1111
$ws.Range("J134").Value = 34
$ws.Range("K134").Value = 5196.3333
$ws.Range("L134").Value = 102
$ws.Range("M134").Value = -2661.3333
$ws.Range("N134").Value = -5172

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1210.3529
$ws.Range("I16").Value = 1277.2307
$ws.Range("K16").Value = 1277.2307
$ws.Range("M16").Value = -990.2307000000001
$ws.Range("H68").Value = 65588.336
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("H71").Value = 65588.336
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("H99").Value = 1770.7273
$ws.Range("I99").Value = 1639.7142
$ws.Range("K99").Value = 1639.7142
$ws.Range("M99").Value = -141.7141999999999
$ws.Range("H113").Value = 1210.3529
$ws.Range("I113").Value = 1277.2307
$ws.Range("K113").Value = 1277.2307
$ws.Range("M113").Value = 892.7692999999999
$ws.Range("H126").Value = 1770.7273
$ws.Range("I126").Value = 1639.7142
$ws.Range("K126").Value = 4919.142599999999
$ws.Range("M126").Value = -2449.142599999999
$ws.Range("H134").Value = 3461.625
$ws.Range("I134").Value = 2525.5715
$ws.Range("K134").Value = 7576.7145
$ws.Range("M134").Value = -5041.7145
$ws.Range("M68").ClearContents()
$ws.Range("M71").ClearContents()

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1096.25
$ws.Range("I5").Value = 1118
$ws.Range("K5").Value = 3354
$ws.Range("M5").Value = -3242
$ws.Range("H49").Value = 20
$ws.Range("I49").Value = 20
$ws.Range("K49").Value = 60
$ws.Range("M49").Value = 96
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H135").Value = 1096.25
$ws.Range("I135").Value = 1118
$ws.Range("K135").Value = 10062
$ws.Range("M135").Value = -7527
$ws.Range("H136").Value = 12000
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("M136").ClearContents()

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 8210.6
$ws.Range("J57").Value = 14999
$ws.Range("L57").Value = 14999
$ws.Range("N57").Value = -16639
$ws.Range("H97").Value = 316.27274
$ws.Range("I97").Value = 283.64706
$ws.Range("J97").Value = 427.2
$ws.Range("K97").Value = 283.64706
$ws.Range("L97").Value = 427.2
$ws.Range("M97").Value = 212.35294
$ws.Range("N97").Value = -1419.2
$ws.Range("H113").Value = 5490.923
$ws.Range("I113").Value = 4561.857
$ws.Range("J113").Value = 6574.8335
$ws.Range("K113").Value = 4561.857
$ws.Range("L113").Value = 6574.8335
$ws.Range("M113").Value = -2391.857
$ws.Range("N113").Value = -10914.8335
$ws.Range("H132").Value = 2141.8
$ws.Range("I132").Value = 2141.8
$ws.Range("K132").Value = 6425.400000000001
$ws.Range("M132").Value = -3895.400000000001

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 640.1429000000001
$ws.Range("I93").Value = 640.1429000000001
$ws.Range("K93").Value = 640.1429000000001
$ws.Range("M93").Value = 607.8570999999999
$ws.Range("H104").Value = 10200
$ws.Range("J104").Value = 10200
$ws.Range("L104").Value = 10200
$ws.Range("N104").Value = -17188
$ws.Range("H136").Value = 3009.5
$ws.Range("I136").Value = 2725.2856
$ws.Range("J136").Value = 4999
$ws.Range("K136").Value = 8175.8568
$ws.Range("L136").Value = 14997
$ws.Range("M136").Value = -5625.8568
$ws.Range("N136").Value = -20097

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 51500
$ws.Range("J64").Value = 51500
$ws.Range("L64").Value = 51500
$ws.Range("N64").Value = -51996
$ws.Range("H67").Value = 51500
$ws.Range("J67").Value = 51500
$ws.Range("L67").Value = 51500
$ws.Range("N67").Value = -53216
$ws.Range("H96").Value = 950
$ws.Range("I96").Value = 950
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 950
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 423
$ws.Range("H107").Value = 355
$ws.Range("I107").Value = 355
$ws.Range("K107").Value = 1065
$ws.Range("M107").Value = 855
$ws.Range("H136").Value = 2777.4194
$ws.Range("I136").Value = 1874.5714
$ws.Range("J136").Value = 4673.4
$ws.Range("K136").Value = 5623.7142
$ws.Range("L136").Value = 14020.2
$ws.Range("M136").Value = -3073.7142
$ws.Range("N136").Value = -19120.2
$ws.Range("N96").ClearContents()

